$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns per data row, in order: F(trial_total), H(category), I(cond_cat),
# K(correct_answer), L(stimulus), M(conceptual), N(perceptual), O(typicality),
# P(n), Q(p_typicality), R(p_conceptual), S(p_perceptual)
$data = @(
    @(269, "kitchens", "target", "j", "stimuli/img_d8xbu.png", 91.36363636363636, 73.18181818181819, 82.27272727272728, 33, 10, 10, 10),
    @(270, "living_rooms", "distractor", "f", "stimuli/img_kq9s9.png", 62.30232558139535, 39.97674418604651, 51.13953488372093, 43, 4, 4, 4),
    @(271, "kitchens", "target", "j", "stimuli/img_yeh72.png", 68.66666666666667, 45.21212121212121, 56.93939393939394, 33, 4, 4, 4),
    @(272, "kitchens", "target", "j", "stimuli/img_aplao.png", 64.0909090909091, 40.75757575757576, 52.42424242424242, 33, 3, 3, 3),
    @(273, "kitchens", "target", "j", "stimuli/img_7wul8.png", 43.03030303030303, 25.54545454545455, 34.28787878787879, 33, 1, 1, 1),
    @(274, "living_rooms", "distractor", "f", "stimuli/img_37hgm.png", 70.95454545454545, 54.77272727272727, 62.86363636363636, 44, 6, 6, 6),
    @(275, "kitchens", "target", "j", "stimuli/img_iyxnj.png", 75.30555555555556, 54.33333333333334, 64.81944444444444, 36, 6, 6, 6),
    @(276, "kitchens", "target", "j", "stimuli/img_a8wvq.png", 86.25925925925925, 66.25925925925925, 76.25925925925925, 27, 10, 10, 10),
    @(277, "kitchens", "target", "j", "stimuli/img_ye5sl.png", 53.2258064516129, 34.45161290322581, 43.83870967741936, 31, 2, 2, 2),
    @(278, "kitchens", "target", "j", "stimuli/img_inqod.png", 70.84848484848484, 50.63636363636363, 60.74242424242424, 33, 5, 5, 5),
    @(279, "kitchens", "target", "j", "stimuli/img_ce9vx.png", 75.9090909090909, 57.12121212121212, 66.51515151515152, 33, 7, 7, 7),
    @(280, "kitchens", "target", "j", "stimuli/img_wyl6z.png", 59.8235294117647, 36.23529411764706, 48.02941176470588, 34, 3, 3, 3),
    @(281, "kitchens", "target", "j", "stimuli/img_t90e2.png", 83.0625, 61.96875, 72.515625, 32, 9, 9, 9),
    @(282, "kitchens", "target", "j", "stimuli/img_60242.png", 78.33333333333333, 57.57575757575758, 67.95454545454545, 33, 7, 7, 7),
    @(283, "kitchens", "target", "j", "stimuli/img_nyv2b.png", 11.91176470588235, 6.852941176470588, 9.382352941176471, 34, 1, 1, 1),
    @(284, "living_rooms", "distractor", "f", "stimuli/img_iudc4.png", 73.625, 52.275, 62.95, 40, 6, 6, 6),
    @(285, "kitchens", "target", "j", "stimuli/img_9mky8.png", 84.32352941176471, 65.17647058823529, 74.75, 34, 9, 9, 9),
    @(286, "living_rooms", "distractor", "f", "stimuli/img_3sw8t.png", 67.4888888888889, 48.51111111111111, 58, 45, 5, 5, 5),
    @(287, "kitchens", "target", "j", "stimuli/img_es7o2.png", 52.48571428571429, 27.54285714285714, 40.01428571428572, 35, 2, 2, 2),
    @(288, "kitchens", "target", "j", "stimuli/img_eatdk.png", 81.40625, 61.375, 71.390625, 32, 8, 8, 8),
    @(289, "kitchens", "target", "j", "stimuli/img_uwv6y.png", 78.88888888888889, 59.30555555555556, 69.09722222222223, 36, 8, 8, 8),
    @(290, "kitchens", "target", "j", "stimuli/img_q1ynd.png", 70.05714285714286, 47.31428571428572, 58.68571428571429, 35, 5, 5, 5),
    @(291, "kitchens", "target", "j", "stimuli/img_79b5l.png", 72.74285714285715, 53.31428571428572, 63.02857142857143, 35, 6, 6, 6),
    @(292, "bedrooms", "distractor", "f", "stimuli/img_n9xll.png", 77.14285714285714, 59.21428571428572, 68.17857142857143, 42, 7, 7, 7),
    @(293, "kitchens", "target", "j", "stimuli/img_cv6mf.png", 66.8, 42.08, 54.44, 25, 4, 4, 4),
    @(294, "living_rooms", "distractor", "f", "stimuli/img_hmmra.png", 54.65853658536585, 34.24390243902439, 44.45121951219512, 41, 3, 3, 3)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 6).Value = $row[0]
    $ws.Cells.Item($r, 8).Value = $row[1]
    $ws.Cells.Item($r, 9).Value = $row[2]
    $ws.Cells.Item($r, 11).Value = $row[3]
    $ws.Cells.Item($r, 12).Value = $row[4]
    $ws.Cells.Item($r, 13).Value = $row[5]
    $ws.Cells.Item($r, 14).Value = $row[6]
    $ws.Cells.Item($r, 15).Value = $row[7]
    $ws.Cells.Item($r, 16).Value = $row[8]
    $ws.Cells.Item($r, 17).Value = $row[9]
    $ws.Cells.Item($r, 18).Value = $row[10]
    $ws.Cells.Item($r, 19).Value = $row[11]
    $r = $r + 1
}
